# Apply the "freshly started new application" edit:
# - Remove Sheet3 entirely
# - Trim Sheet2 down to header row + a single new task row
# - Make Sheet2 the active/selected sheet

$wb = $excel.ActiveWorkbook

# --- Update Sheet2 contents ---
$ws2 = $wb.Worksheets.Item("Sheet2")

# Clear out old rows 3-8 (keep header row 1)
$ws2.Rows("3:8").Delete() | Out-Null

# Set the new task text in A3
$ws2.Range("A3").Value = "not throwing my exception throwing null pointer"

# Update selection to A3
$ws2.Range("A3").Select() | Out-Null

# --- Remove Sheet3 ---
$excel.DisplayAlerts = $false
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Delete() | Out-Null
$excel.DisplayAlerts = $true

# --- Make Sheet2 the active sheet ---
$ws2.Activate() | Out-Null
